$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 1434.8334
$ws.Range("I31").Value = 121.8
$ws.Range("J31").Value = 8000
$ws.Range("K31").Value = 365.4
$ws.Range("L31").Value = 24000
$ws.Range("M31").Value = -135.4
$ws.Range("N31").Value = -24460

$ws.Range("H70").Value = 1194.375
$ws.Range("I70").Value = 999
$ws.Range("J70").Value = 1274.8235
$ws.Range("K70").Value = 2997
$ws.Range("L70").Value = 3824.4705
$ws.Range("M70").Value = -2727
$ws.Range("N70").Value = -4364.470499999999

$ws.Range("H73").Value = 1194.375
$ws.Range("I73").Value = 999
$ws.Range("J73").Value = 1274.8235
$ws.Range("K73").Value = 2997
$ws.Range("L73").Value = 3824.4705
$ws.Range("M73").Value = -2061
$ws.Range("N73").Value = -5696.470499999999

$ws.Range("H80").Value = 1857.0741
$ws.Range("I80").Value = 1317
$ws.Range("J80").Value = 2289.1333
$ws.Range("K80").Value = 3951
$ws.Range("L80").Value = 6867.3999
$ws.Range("M80").Value = -2953
$ws.Range("N80").Value = -8863.3999

$ws.Range("H82").Value = 821
$ws.Range("I82").Value = 821
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 2463
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -2057
$ws.Range("N82").ClearContents()

$ws.Range("H83").Value = 1857.0741
$ws.Range("I83").Value = 1317
$ws.Range("J83").Value = 2289.1333
$ws.Range("K83").Value = 11853
$ws.Range("L83").Value = 20602.1997
$ws.Range("M83").Value = -6861
$ws.Range("N83").Value = -30586.1997

$ws.Range("H85").Value = 821
$ws.Range("I85").Value = 821
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 2463
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -1059
$ws.Range("N85").ClearContents()

$ws.Range("H86").Value = 5229.0713
$ws.Range("I86").Value = 6334.3335
$ws.Range("J86").Value = 4927.636
$ws.Range("K86").Value = 6334.3335
$ws.Range("L86").Value = 4927.636
$ws.Range("M86").Value = -5211.3335
$ws.Range("N86").Value = -7173.636

$ws.Range("H88").Value = 21946.227
$ws.Range("I88").Value = 48144.75
$ws.Range("J88").Value = 5399.7896
$ws.Range("K88").Value = 48144.75
$ws.Range("L88").Value = 5399.7896
$ws.Range("M88").Value = -47738.75
$ws.Range("N88").Value = -6211.7896

$ws.Range("H89").Value = 5229.0713
$ws.Range("I89").Value = 6334.3335
$ws.Range("J89").Value = 4927.636
$ws.Range("K89").Value = 31671.6675
$ws.Range("L89").Value = 24638.18
$ws.Range("M89").Value = -26055.6675
$ws.Range("N89").Value = -35870.18

$ws.Range("H91").Value = 21946.227
$ws.Range("I91").Value = 48144.75
$ws.Range("J91").Value = 5399.7896
$ws.Range("K91").Value = 48144.75
$ws.Range("L91").Value = 5399.7896
$ws.Range("M91").Value = -46740.75
$ws.Range("N91").Value = -8207.7896

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1600
$ws.Range("I88").Value = 1600
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 1600
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -1194
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 1600
$ws.Range("I91").Value = 1600
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 1600
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -196
$ws.Range("N91").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 20140
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 20140
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 20140
$ws.Range("N26").Value = -20714

$ws.Range("H62").Value = 2281.3157
$ws.Range("I62").Value = 2237.4167
$ws.Range("J62").Value = 2356.5715
$ws.Range("K62").Value = 2237.4167
$ws.Range("L62").Value = 2356.5715
$ws.Range("M62").Value = -1613.4167
$ws.Range("N62").Value = -3604.5715

$ws.Range("H65").Value = 2281.3157
$ws.Range("I65").Value = 2237.4167
$ws.Range("J65").Value = 2356.5715
$ws.Range("K65").Value = 11187.0835
$ws.Range("L65").Value = 11782.8575
$ws.Range("M65").Value = -8067.083500000001
$ws.Range("N65").Value = -18022.8575

$ws.Range("H70").Value = 28833.334
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 28833.334
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 28833.334
$ws.Range("N70").Value = -29463.334

$ws.Range("H73").Value = 28833.334
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 28833.334
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 28833.334
$ws.Range("N73").Value = -31017.334

$ws.Range("H75").Value = 24619.25
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 24619.25
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 24619.25
$ws.Range("N75").Value = -26615.25

$ws.Range("H78").Value = 24619.25
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 24619.25
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 73857.75
$ws.Range("N78").Value = -83841.75

$ws.Range("H86").Value = 11054.714
$ws.Range("I86").Value = 4725.75
$ws.Range("J86").Value = 19493.334
$ws.Range("K86").Value = 4725.75
$ws.Range("L86").Value = 19493.334
$ws.Range("M86").Value = -3602.75
$ws.Range("N86").Value = -21739.334

$ws.Range("H89").Value = 11054.714
$ws.Range("I89").Value = 4725.75
$ws.Range("J89").Value = 19493.334
$ws.Range("K89").Value = 23628.75
$ws.Range("L89").Value = 97466.67
$ws.Range("M89").Value = -18012.75
$ws.Range("N89").Value = -108698.67

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1874.091
$ws.Range("I39").Value = 1073.5714
$ws.Range("J39").Value = 3275
$ws.Range("K39").Value = 3220.7142
$ws.Range("L39").Value = 9825
$ws.Range("M39").Value = -2926.7142
$ws.Range("N39").Value = -10413

$ws.Range("H49").Value = 1201.8

$ws.Range("H113").Value = 947592.25
$ws.Range("I113").Value = 2331625.2
$ws.Range("J113").Value = 622.2632
$ws.Range("K113").Value = 6994875.600000001
$ws.Range("L113").Value = 1866.7896
$ws.Range("M113").Value = -6992705.600000001
$ws.Range("N113").Value = -6206.7896

$ws.Range("H131").Value = 918.77
$ws.Range("I131").Value = 312
$ws.Range("J131").Value = 944.05206
$ws.Range("K131").Value = 936
$ws.Range("L131").Value = 2832.15618
$ws.Range("M131").Value = 4104
$ws.Range("N131").Value = -12912.15618

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2153.8462
$ws.Range("I80").Value = 2133.3333
$ws.Range("J80").Value = 2400
$ws.Range("K80").Value = 2133.3333
$ws.Range("L80").Value = 2400
$ws.Range("M80").Value = -1135.3333
$ws.Range("N80").Value = -4396

$ws.Range("H83").Value = 2153.8462
$ws.Range("I83").Value = 2133.3333
$ws.Range("J83").Value = 2400
$ws.Range("K83").Value = 10666.6665
$ws.Range("L83").Value = 12000
$ws.Range("M83").Value = -5674.666499999999
$ws.Range("N83").Value = -21984

$ws.Range("H107").Value = 774.1429000000001
$ws.Range("I107").Value = 717.5
$ws.Range("J107").Value = 915.75
$ws.Range("K107").Value = 717.5
$ws.Range("L107").Value = 915.75
$ws.Range("M107").Value = 1202.5
$ws.Range("N107").Value = -4755.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 10125
$ws.Range("I68").Value = 18016.666
$ws.Range("J68").Value = 2233.3333
$ws.Range("K68").Value = 18016.666
$ws.Range("L68").Value = 2233.3333
$ws.Range("M68").Value = -17267.666
$ws.Range("N68").Value = -3731.3333

$ws.Range("H71").Value = 10125
$ws.Range("I71").Value = 18016.666
$ws.Range("J71").Value = 2233.3333
$ws.Range("K71").Value = 90083.33
$ws.Range("L71").Value = 11166.6665
$ws.Range("M71").Value = -86339.33
$ws.Range("N71").Value = -18654.6665

$ws.Range("H100").Value = 1184.1875
$ws.Range("I100").Value = 1003
$ws.Range("J100").Value = 1727.75
$ws.Range("K100").Value = 1003
$ws.Range("L100").Value = 1727.75
$ws.Range("M100").Value = -462

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1618.625
$ws.Range("I81").Value = 1333
$ws.Range("J81").Value = 1790
$ws.Range("K81").Value = 2666
$ws.Range("L81").Value = 3580
$ws.Range("M81").Value = -1605
$ws.Range("N81").Value = -5702

$ws.Range("H84").Value = 1618.625
$ws.Range("I84").Value = 1333
$ws.Range("J84").Value = 1790
$ws.Range("K84").Value = 13330
$ws.Range("L84").Value = 17900
$ws.Range("M84").Value = -8026
$ws.Range("N84").Value = -28508

$ws.Range("H86").Value = 27000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 27000
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 27000
$ws.Range("N86").Value = -29246

$ws.Range("H89").Value = 27000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 27000
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 135000
$ws.Range("N89").Value = -146232

$ws.Range("H136").Value = 11440.533
$ws.Range("I136").Value = 13427
$ws.Range("J136").Value = 10718.182
$ws.Range("K136").Value = 40281
$ws.Range("L136").Value = 32154.546
$ws.Range("M136").Value = -37731
$ws.Range("N136").Value = -37254.546
